$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.914.43'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.628.62'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'211.64"
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'23.35"
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = "'0.0880"
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = '1.859.20'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '1.616.82'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('E15').Value = '  -2.10%  '
$ws.Range('D16').Value = "'65.62"
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '27.893.95'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = "'230.64"
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = "'7.68"
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').Value = '0.0₃0724'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = "'0.999"
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').Value = "'10.25"
$ws.Range('E23').Value = '  -5.14%  '
$ws.Range('D24').Value = "'2.05"
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').Value = "'154.66"
$ws.Range('E25').Value = '  +2.05%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = "'15.56"
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -0.82%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = '1.401.08'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('E36').Value = '  +9.38%  '
$ws.Range('D37').Value = "'2.34"
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').Value = "'0.863"
$ws.Range('E40').Value = '  -2.74%  '
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').Value = "'0.999"
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').Value = "'1.85"
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').Value = "'66.16"
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('D47').Value = '1.768.65'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('D48').Value = "'88.13"
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('E51').Value = '  -0.31%  '
